# update scripts wuth new tpm
# Refresh the NATMI ligand-receptor (Ptn -> Plxnb2) metrics with the
# values recomputed from the new TPM-normalised expression matrix.
# Ligand/receptor expression stats and every derived "Edge ..." column
# change as a consequence; cluster/ligand/receptor labels (cols A-D) and
# the expressing-cell counts/detection-rate columns K/L stay the same.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.7169949999999999
$ws.Cells.Item(2, 8).Value = 2.150985
$ws.Cells.Item(2, 9).Value = 0.003651663653539308
$ws.Cells.Item(2, 10).Value = 0.003651663653539308
$ws.Cells.Item(2, 13).Value = 10.82167433333333
$ws.Cells.Item(2, 14).Value = 32.465023
$ws.Cells.Item(2, 15).Value = 0.09133543757015983
$ws.Cells.Item(2, 16).Value = 0.09133543757015983
$ws.Cells.Item(2, 17).Value = 7.759086388628333
$ws.Cells.Item(2, 18).Value = 69.831777497655
$ws.Cells.Item(2, 19).Value = 0.0003335262976550612
$ws.Cells.Item(2, 20).Value = 0.0003335262976550612

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.7169949999999999
$ws.Cells.Item(3, 8).Value = 2.150985
$ws.Cells.Item(3, 9).Value = 0.003651663653539308
$ws.Cells.Item(3, 10).Value = 0.003651663653539308
$ws.Cells.Item(3, 13).Value = 36.14140700000001
$ws.Cells.Item(3, 15).Value = 0.3050351656377608
$ws.Cells.Item(3, 16).Value = 0.3050351656377608
$ws.Cells.Item(3, 17).Value = 25.91320811196501
$ws.Cells.Item(3, 18).Value = 233.218873007685
$ws.Cells.Item(3, 19).Value = 0.001113885827410754
$ws.Cells.Item(3, 20).Value = 0.001113885827410754

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.7169949999999999
$ws.Cells.Item(4, 8).Value = 2.150985
$ws.Cells.Item(4, 9).Value = 0.003651663653539308
$ws.Cells.Item(4, 10).Value = 0.003651663653539308
$ws.Cells.Item(4, 13).Value = 26.40107466666666
$ws.Cells.Item(4, 14).Value = 79.20322399999999
$ws.Cells.Item(4, 15).Value = 0.2228263051286729
$ws.Cells.Item(4, 16).Value = 0.2228263051286729
$ws.Cells.Item(4, 17).Value = 18.92943853062666
$ws.Cells.Item(4, 18).Value = 170.36494677564
$ws.Cells.Item(4, 19).Value = 0.0008136867194908343
$ws.Cells.Item(4, 20).Value = 0.0008136867194908344

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.7169949999999999
$ws.Cells.Item(5, 8).Value = 2.150985
$ws.Cells.Item(5, 9).Value = 0.003651663653539308
$ws.Cells.Item(5, 10).Value = 0.003651663653539308
$ws.Cells.Item(5, 13).Value = 45.11859966666666
$ws.Cells.Item(5, 14).Value = 135.355799
$ws.Cells.Item(5, 15).Value = 0.3808030916634065
$ws.Cells.Item(5, 16).Value = 0.3808030916634065
$ws.Cells.Item(5, 17).Value = 32.34981036800166
$ws.Cells.Item(5, 18).Value = 291.148293312015
$ws.Cells.Item(5, 19).Value = 0.001390564808982659
$ws.Cells.Item(5, 20).Value = 0.001390564808982659

# Row 6
$ws.Cells.Item(6, 9).Value = 0.9751961860217362
$ws.Cells.Item(6, 10).Value = 0.9751961860217361
$ws.Cells.Item(6, 13).Value = 10.82167433333333
$ws.Cells.Item(6, 14).Value = 32.465023
$ws.Cells.Item(6, 15).Value = 0.09133543757015983
$ws.Cells.Item(6, 16).Value = 0.09133543757015983
$ws.Cells.Item(6, 17).Value = 2072.10525697505
$ws.Cells.Item(6, 18).Value = 18648.94731277545
$ws.Cells.Item(6, 19).Value = 0.08906997036704625
$ws.Cells.Item(6, 20).Value = 0.08906997036704624

# Row 7
$ws.Cells.Item(7, 9).Value = 0.9751961860217362
$ws.Cells.Item(7, 10).Value = 0.9751961860217361
$ws.Cells.Item(7, 13).Value = 36.14140700000001
$ws.Cells.Item(7, 15).Value = 0.3050351656377608
$ws.Cells.Item(7, 16).Value = 0.3050351656377608
$ws.Cells.Item(7, 17).Value = 6920.259946143411
$ws.Cells.Item(7, 19).Value = 0.2974691301324529
$ws.Cells.Item(7, 20).Value = 0.2974691301324529

# Row 8
$ws.Cells.Item(8, 9).Value = 0.9751961860217362
$ws.Cells.Item(8, 10).Value = 0.9751961860217361
$ws.Cells.Item(8, 13).Value = 26.40107466666666
$ws.Cells.Item(8, 14).Value = 79.20322399999999
$ws.Cells.Item(8, 15).Value = 0.2228263051286729
$ws.Cells.Item(8, 16).Value = 0.2228263051286729
$ws.Cells.Item(8, 17).Value = 5055.207163098957
$ws.Cells.Item(8, 18).Value = 45496.86446789061
$ws.Cells.Item(8, 19).Value = 0.2172993629067974
$ws.Cells.Item(8, 20).Value = 0.2172993629067974

# Row 9
$ws.Cells.Item(9, 9).Value = 0.9751961860217362
$ws.Cells.Item(9, 10).Value = 0.9751961860217361
$ws.Cells.Item(9, 13).Value = 45.11859966666666
$ws.Cells.Item(9, 14).Value = 135.355799
$ws.Cells.Item(9, 15).Value = 0.3808030916634065
$ws.Cells.Item(9, 16).Value = 0.3808030916634065
$ws.Cells.Item(9, 17).Value = 8639.188786958755
$ws.Cells.Item(9, 18).Value = 77752.69908262881
$ws.Cells.Item(9, 19).Value = 0.3713577226154396
$ws.Cells.Item(9, 20).Value = 0.3713577226154396

# Row 10
$ws.Cells.Item(10, 7).Value = 4.138615666666666
$ws.Cells.Item(10, 8).Value = 12.415847
$ws.Cells.Item(10, 9).Value = 0.02107801645190694
$ws.Cells.Item(10, 10).Value = 0.02107801645190694
$ws.Cells.Item(10, 13).Value = 10.82167433333333
$ws.Cells.Item(10, 14).Value = 32.465023
$ws.Cells.Item(10, 15).Value = 0.09133543757015983
$ws.Cells.Item(10, 16).Value = 0.09133543757015983
$ws.Cells.Item(10, 17).Value = 44.78675093549789
$ws.Cells.Item(10, 18).Value = 403.080758419481
$ws.Cells.Item(10, 19).Value = 0.001925169855745949
$ws.Cells.Item(10, 20).Value = 0.001925169855745948

# Row 11
$ws.Cells.Item(11, 7).Value = 4.138615666666666
$ws.Cells.Item(11, 8).Value = 12.415847
$ws.Cells.Item(11, 9).Value = 0.02107801645190694
$ws.Cells.Item(11, 10).Value = 0.02107801645190694
$ws.Cells.Item(11, 13).Value = 36.14140700000001
$ws.Cells.Item(11, 15).Value = 0.3050351656377608
$ws.Cells.Item(11, 16).Value = 0.3050351656377608
$ws.Cells.Item(11, 17).Value = 149.5753932255764
$ws.Cells.Item(11, 18).Value = 1346.178539030187
$ws.Cells.Item(11, 19).Value = 0.006429536239722882
$ws.Cells.Item(11, 20).Value = 0.006429536239722881

# Row 12
$ws.Cells.Item(12, 7).Value = 4.138615666666666
$ws.Cells.Item(12, 8).Value = 12.415847
$ws.Cells.Item(12, 9).Value = 0.02107801645190694
$ws.Cells.Item(12, 10).Value = 0.02107801645190694
$ws.Cells.Item(12, 13).Value = 26.40107466666666
$ws.Cells.Item(12, 14).Value = 79.20322399999999
$ws.Cells.Item(12, 15).Value = 0.2228263051286729
$ws.Cells.Item(12, 16).Value = 0.2228263051286729
$ws.Cells.Item(12, 17).Value = 109.2639012323031
$ws.Cells.Item(12, 18).Value = 983.3751110907278
$ws.Cells.Item(12, 19).Value = 0.004696736525419804
$ws.Cells.Item(12, 20).Value = 0.004696736525419803

# Row 13
$ws.Cells.Item(13, 7).Value = 4.138615666666666
$ws.Cells.Item(13, 8).Value = 12.415847
$ws.Cells.Item(13, 9).Value = 0.02107801645190694
$ws.Cells.Item(13, 10).Value = 0.02107801645190694
$ws.Cells.Item(13, 13).Value = 45.11859966666666
$ws.Cells.Item(13, 14).Value = 135.355799
$ws.Cells.Item(13, 15).Value = 0.3808030916634065
$ws.Cells.Item(13, 16).Value = 0.3808030916634065
$ws.Cells.Item(13, 17).Value = 186.7285434385281
$ws.Cells.Item(13, 18).Value = 1680.556890946753
$ws.Cells.Item(13, 19).Value = 0.00802657383101831
$ws.Cells.Item(13, 20).Value = 0.00802657383101831

# Row 14
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.014556
$ws.Cells.Item(14, 8).Value = 0.043668
$ws.Cells.Item(14, 9).Value = 0.00007413387281768795
$ws.Cells.Item(14, 10).Value = 0.00007413387281768795
$ws.Cells.Item(14, 13).Value = 10.82167433333333
$ws.Cells.Item(14, 14).Value = 32.465023
$ws.Cells.Item(14, 15).Value = 0.09133543757015983
$ws.Cells.Item(14, 16).Value = 0.09133543757015983
$ws.Cells.Item(14, 17).Value = 0.157520291596
$ws.Cells.Item(14, 18).Value = 1.417682624364
$ws.Cells.Item(14, 19).Value = 0.000006771049712574106
$ws.Cells.Item(14, 20).Value = 0.000006771049712574106

# Row 15
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.014556
$ws.Cells.Item(15, 8).Value = 0.043668
$ws.Cells.Item(15, 9).Value = 0.00007413387281768795
$ws.Cells.Item(15, 10).Value = 0.00007413387281768795
$ws.Cells.Item(15, 13).Value = 36.14140700000001
$ws.Cells.Item(15, 15).Value = 0.3050351656377608
$ws.Cells.Item(15, 16).Value = 0.3050351656377608
$ws.Cells.Item(15, 17).Value = 0.5260743202920001
$ws.Cells.Item(15, 18).Value = 4.734668882628001
$ws.Cells.Item(15, 19).Value = 0.00002261343817431213
$ws.Cells.Item(15, 20).Value = 0.00002261343817431213

# Row 16
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.014556
$ws.Cells.Item(16, 8).Value = 0.043668
$ws.Cells.Item(16, 9).Value = 0.00007413387281768795
$ws.Cells.Item(16, 10).Value = 0.00007413387281768795
$ws.Cells.Item(16, 13).Value = 26.40107466666666
$ws.Cells.Item(16, 14).Value = 79.20322399999999
$ws.Cells.Item(16, 15).Value = 0.2228263051286729
$ws.Cells.Item(16, 16).Value = 0.2228263051286729
$ws.Cells.Item(16, 17).Value = 0.384294042848
$ws.Cells.Item(16, 18).Value = 3.458646385632
$ws.Cells.Item(16, 19).Value = 0.00001651897696484436
$ws.Cells.Item(16, 20).Value = 0.00001651897696484436

# Row 17
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.014556
$ws.Cells.Item(17, 8).Value = 0.043668
$ws.Cells.Item(17, 9).Value = 0.00007413387281768795
$ws.Cells.Item(17, 10).Value = 0.00007413387281768795
$ws.Cells.Item(17, 13).Value = 45.11859966666666
$ws.Cells.Item(17, 14).Value = 135.355799
$ws.Cells.Item(17, 15).Value = 0.3808030916634065
$ws.Cells.Item(17, 16).Value = 0.3808030916634065
$ws.Cells.Item(17, 17).Value = 0.6567463367479999
$ws.Cells.Item(17, 18).Value = 5.910717030731999
$ws.Cells.Item(17, 19).Value = 0.00002823040796595734
$ws.Cells.Item(17, 20).Value = 0.00002823040796595734
